$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.08160078149282
$ws.Range("D2").Value = 1.081037149221241
$ws.Range("E2").Value = 1.08407961756302
$ws.Range("F2").Value = 1.093804605640814
$ws.Range("I2").Value = 1.053819493488765
$ws.Range("J2").Value = 1.086474558061409
$ws.Range("K2").Value = 1.083709018803914
$ws.Range("L2").Value = 1.08674356793931
$ws.Range("M2").Value = 1.096443561597217
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.083401393453988
$ws.Range("D3").Value = 1.082480895611814
$ws.Range("E3").Value = 1.085685934394403
$ws.Range("F3").Value = 1.095454131505379
$ws.Range("I3").Value = 1.05432756573479
$ws.Range("J3").Value = 1.087933031224308
$ws.Range("K3").Value = 1.084969620202296
$ws.Range("L3").Value = 1.088166927195285
$ws.Range("M3").Value = 1.097911859235938
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.084563598509203
$ws.Range("D4").Value = 1.083412368181115
$ws.Range("E4").Value = 1.086722425092756
$ws.Range("F4").Value = 1.096518839316773
$ws.Range("I4").Value = 1.054653478116228
$ws.Range("J4").Value = 1.088873509858958
$ws.Range("K4").Value = 1.085782062731296
$ws.Range("L4").Value = 1.089084541136353
$ws.Range("M4").Value = 1.098858813018984
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.085051508273729
$ws.Range("D5").Value = 1.083803317874501
$ws.Range("E5").Value = 1.087157484349062
$ws.Range("F5").Value = 1.096965822618769
$ws.Range("I5").Value = 1.054789815877302
$ws.Range("J5").Value = 1.089268121634237
$ws.Range("K5").Value = 1.086122846318214
$ws.Range("L5").Value = 1.08946950619243
$ws.Range("M5").Value = 1.099256175016162
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.085133390976264
$ws.Range("D6").Value = 1.083868922775308
$ws.Range("E6").Value = 1.087230493149405
$ws.Range("F6").Value = 1.097040837123897
$ws.Range("I6").Value = 1.054812668106053
$ws.Range("J6").Value = 1.089334334173259
$ws.Range("K6").Value = 1.08618002070799
$ws.Range("L6").Value = 1.089534096987207
$ws.Range("M6").Value = 1.099322850969975
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.084570120638652
$ws.Range("D7").Value = 1.083417594574163
$ws.Range("E7").Value = 1.086728241034599
$ws.Range("F7").Value = 1.096524814346609
$ws.Range("I7").Value = 1.054655302516434
$ws.Range("J7").Value = 1.088878785670854
$ws.Range("K7").Value = 1.085786619297922
$ws.Range("L7").Value = 1.089089688183513
$ws.Range("M7").Value = 1.09886412546944
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.082209919042333
$ws.Range("D8").Value = 1.081525641966027
$ws.Range("E8").Value = 1.084623089379467
$ws.Range("F8").Value = 1.094362626622462
$ws.Range("I8").Value = 1.053991790928023
$ws.Range("J8").Value = 1.08696813681553
$ws.Range("K8").Value = 1.084135724720966
$ws.Range("L8").Value = 1.087225309518891
$ws.Range("M8").Value = 1.096940436017959
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.078027936590262
$ws.Range("D9").Value = 1.078170360271453
$ws.Range("E9").Value = 1.08089068606773
$ws.Range("F9").Value = 1.090531702637809
$ws.Range("I9").Value = 1.052800582562804
$ws.Range("J9").Value = 1.083575850690425
$ws.Range("K9").Value = 1.081201243022051
$ws.Range("L9").Value = 1.083913457955661
$ws.Range("M9").Value = 1.093526083557882
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.075223485393484
$ws.Range("D10").Value = 1.075918363237348
$ws.Range("E10").Value = 1.078386181003755
$ws.Range("F10").Value = 1.087962862545596
$ws.Range("I10").Value = 1.051991314366398
$ws.Range("J10").Value = 1.081296381485388
$ws.Range("K10").Value = 1.079227149437758
$ws.Range("L10").Value = 1.081686873956055
$ws.Range("M10").Value = 1.091232529364708
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.0740049996537
$ws.Range("D11").Value = 1.074939465857998
$ws.Range("E11").Value = 1.077297657323718
$ws.Range("F11").Value = 1.086846802646728
$ws.Range("I11").Value = 1.051637230362952
$ws.Range("J11").Value = 1.080304907380444
$ws.Range("K11").Value = 1.078367973120841
$ws.Range("L11").Value = 1.080718125116163
$ws.Range("M11").Value = 1.090235107456671
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.073551758549812
$ws.Range("D12").Value = 1.074575279288094
$ws.Range("E12").Value = 1.07689270488419
$ws.Range("F12").Value = 1.08643167050273
$ws.Range("I12").Value = 1.051505150679073
$ws.Range("J12").Value = 1.079935945532258
$ws.Range("K12").Value = 1.078048165017124
$ws.Range("L12").Value = 1.080357578283701
$ws.Range("M12").Value = 1.089863959261153
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.073649009634808
$ws.Range("D13").Value = 1.074653425045813
$ws.Range("E13").Value = 1.076979597169938
$ws.Range("F13").Value = 1.086520744172836
$ws.Range("I13").Value = 1.051533507534271
$ws.Range("J13").Value = 1.080015120319566
$ws.Range("K13").Value = 1.078116795559882
$ws.Range("L13").Value = 1.080434949204683
$ws.Range("M13").Value = 1.089943601985933
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.073967547796228
$ws.Range("D14").Value = 1.074909374008582
$ws.Range("E14").Value = 1.077264196702736
$ws.Range("F14").Value = 1.086812499578174
$ws.Range("I14").Value = 1.051626324017541
$ws.Range("J14").Value = 1.080274422929058
$ws.Range("K14").Value = 1.0783415514663
$ws.Range("L14").Value = 1.080688336786009
$ws.Range("M14").Value = 1.090204441811984
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.074163724118568
$ws.Range("D15").Value = 1.075066995255758
$ws.Range("E15").Value = 1.077439464372824
$ws.Range("F15").Value = 1.086992182646816
$ws.Range("I15").Value = 1.051683437299222
$ws.Range("J15").Value = 1.080434096682478
$ws.Range("K15").Value = 1.078479941573018
$ws.Range("L15").Value = 1.080844362630194
$ws.Range("M15").Value = 1.090365065740904
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.075304263416283
$ws.Range("D16").Value = 1.075983248821612
$ws.Range("E16").Value = 1.078458335837358
$ws.Range("F16").Value = 1.088036851653783
$ws.Range("I16").Value = 1.052014735981283
$ws.Range("J16").Value = 1.081362087276772
$ws.Range("K16").Value = 1.079284076685044
$ws.Range("L16").Value = 1.081751067875711
$ws.Range("M16").Value = 1.091298633095226
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.076018571684033
$ws.Range("D17").Value = 1.076556971596372
$ws.Range("E17").Value = 1.079096349129723
$ws.Range("F17").Value = 1.088691133222334
$ws.Range("I17").Value = 1.052221565020735
$ws.Range("J17").Value = 1.081942989049498
$ws.Range("K17").Value = 1.079787306986383
$ws.Range("L17").Value = 1.082318571251627
$ws.Range("M17").Value = 1.091883073636159
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.076434817487739
$ws.Range("D18").Value = 1.076891251829354
$ws.Range("E18").Value = 1.079468101246404
$ws.Range("F18").Value = 1.089072405333712
$ws.Range("I18").Value = 1.052341851666774
$ws.Range("J18").Value = 1.082281391049514
$ws.Range("K18").Value = 1.080080410953791
$ws.Range("L18").Value = 1.082649141365639
$ws.Range("M18").Value = 1.09222355430199
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.076576679601371
$ws.Range("D19").Value = 1.077005171688493
$ws.Range("E19").Value = 1.079594793426022
$ws.Range("F19").Value = 1.089202348834735
$ws.Range("I19").Value = 1.052382806582875
$ws.Range("J19").Value = 1.082396705236334
$ws.Range("K19").Value = 1.080180280704602
$ws.Range("L19").Value = 1.082761782226344
$ws.Range("M19").Value = 1.092339579723146
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.075941974485984
$ws.Range("D20").Value = 1.076495454177193
$ws.Range("E20").Value = 1.07902793682771
$ws.Range("F20").Value = 1.088620972281899
$ws.Range("I20").Value = 1.052199410805093
$ws.Range("J20").Value = 1.081880708189333
$ws.Range("K20").Value = 1.079733358859686
$ws.Range("L20").Value = 1.082257729652947
$ws.Range("M20").Value = 1.091820411591751
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.073873764104524
$ws.Range("D21").Value = 1.074834019559707
$ws.Range("E21").Value = 1.077180406658325
$ws.Range("F21").Value = 1.086726600967909
$ws.Range("I21").Value = 1.051599007307738
$ws.Range("J21").Value = 1.080198083789324
$ws.Range("K21").Value = 1.078275385107802
$ws.Range("L21").Value = 1.080613740196442
$ws.Range("M21").Value = 1.090127649365387
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.072569678383193
$ws.Range("D22").Value = 1.07378604382088
$ws.Range("E22").Value = 1.076015159179035
$ws.Range("F22").Value = 1.085532183154222
$ws.Range("I22").Value = 1.051218283242732
$ws.Range("J22").Value = 1.07913618438741
$ws.Range("K22").Value = 1.077354805663339
$ws.Range("L22").Value = 1.079575980788082
$ws.Range("M22").Value = 1.089059508701493
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.073261357982713
$ws.Range("D23").Value = 1.074341919493559
$ws.Range("E23").Value = 1.076633228852254
$ws.Range("F23").Value = 1.086165689912124
$ws.Range("I23").Value = 1.051420420275434
$ws.Range("J23").Value = 1.079699498675841
$ws.Range("K23").Value = 1.077843195885234
$ws.Range("L23").Value = 1.08012651233938
$ws.Range("M23").Value = 1.089626118935686
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.075976586692565
$ws.Range("D24").Value = 1.076523252371056
$ws.Range("E24").Value = 1.079058850607808
$ws.Range("F24").Value = 1.08865267610072
$ws.Range("I24").Value = 1.052209422432922
$ws.Range("J24").Value = 1.081908851552005
$ws.Range("K24").Value = 1.079757736998993
$ws.Range("L24").Value = 1.082285222727292
$ws.Range("M24").Value = 1.091848727150196
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.079111907392002
$ws.Range("D25").Value = 1.079040395140755
$ws.Range("E25").Value = 1.081858399309327
$ws.Range("F25").Value = 1.091524649524353
$ws.Range("I25").Value = 1.053111180914999
$ws.Range("J25").Value = 1.084455940652284
$ws.Range("K25").Value = 1.081962957700428
$ws.Range("L25").Value = 1.084772883817761
$ws.Range("M25").Value = 1.094411768774962

Write-Host "Updated cells"